# "Removed notes from slides"
#
# The commit clears the speaker-notes body placeholder on slides
# 1,2,3,4,5,6,8 entirely, and rewrites the notes on slide 7 with new text
# (replacing the long "Cancelled, Finished ..." feature-value dump with a
# short paragraph about misleading feature importances). Slides 9 and 10
# already have empty notes and are left untouched.

$p = $ppt.ActivePresentation

# Slide 1 notes: "Hi I'm Michael Engeling. ..."
$s1 = $p.Slides.Item(1)
$s1.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = ""

# Slide 2 notes: "Now a little background on Silvercar, ..."
$s2 = $p.Slides.Item(2)
$s2.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = ""

# Slide 3 notes: "I set out to create a model ..." (+ 4 more bullet notes)
$s3 = $p.Slides.Item(3)
$s3.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = ""

# Slide 4 notes: "I used most of the standard data science tools, ..."
$s4 = $p.Slides.Item(4)
$s4.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = ""

# Slide 5 notes: "As for the data at my disposal, ..." (+ 1 more note)
$s5 = $p.Slides.Item(5)
$s5.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = ""

# Slide 6 notes: "Started with 75% accurate model ..." (+ 5 more notes)
$s6 = $p.Slides.Item(6)
$s6.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = ""

# Slide 7 notes: replace the "Cancelled, Finished" / feature-value dump
# with the new write-up about misleading feature importances.
$s7 = $p.Slides.Item(7)
$newNote = "Gradient boosting classifiers have misleading feature importances when categorical and numerical features are mixed. In this case, the numerical features are weighted more heavily" + "even though all of the features were standardized. Feature importances have to be evaluated separately for the numerical and categorical features."
$s7.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = $newNote

# Slide 8 notes: "Currently only using one user variable, ..."
$s8 = $p.Slides.Item(8)
$s8.NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = ""
